# Apply test5 improvements:
# 1. Fix typo "locacalizacion" -> "localizacion" in B1 header
# 2. Clear the email value in C2 (keep its hyperlink-like style)
# 3. Change the active selection to cell C2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix header typo
$ws.Range("B1").Value = "localizacion"

# 2. Clear C2 content but keep formatting/style
$ws.Range("C2").ClearContents()

# 3. Update selection to C2
$ws.Range("C2").Select()
